# Edit: "nhan xet bai tap ngay 11"
# 1) Remove the stray _GoBack bookmark after "Canh hinh vuong = 12" (it is re-added later,
#    at the end of the "x^2" exercise paragraph).
# 2) Rewrite the x*2/x*4/x*2n exercise paragraph: split the caret text into separate runs,
#    change "*" to "^", and append a note explaining the caret, then re-add the _GoBack bookmark.
# 3) Add w:lang="vi-VN" throughout the "Hay thuc hien tinh tong..." block (4 paragraphs).

$d = $word.ActiveDocument

function Get-ParagraphRangeByText($doc, $needle) {
    $full = $doc.Content.Text
    $idx = $full.IndexOf($needle)
    if ($idx -lt 0) {
        throw "text not found: $needle"
    }
    $para = $doc.Range($idx, $idx + 1).Paragraphs(1)
    return $para.Range
}

# --- Step 1: drop the old bookmark around "Canh hinh vuong = 12m" ---
$bookmarkParaRange = Get-ParagraphRangeByText $d "Cạnh hình vuông"
$bookmarkParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D85F74" w:rsidRPr="00D85F74" w:rsidRDefault="00D85F74" w:rsidP="00D85F74"><w:pPr><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Cạnh hình vuông</w:t></w:r><w:r w:rsidRPr="00D85F74"><w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> = 1</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>2</w:t></w:r><w:r w:rsidRPr="00D85F74"><w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>m</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$bookmarkParaRange.InsertXML($bookmarkParaXml)

# --- Step 2: rewrite the "S(n) = x*2 + x*4 + ... + x*2n" paragraph ---
$tongParaRange = Get-ParagraphRangeByText $d "x*2"
$tongParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0068593B" w:rsidRPr="00C5583D" w:rsidRDefault="001151B4" w:rsidP="001151B4"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="47"/></w:numPr><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r w:rsidRPr="001151B4"><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:lastRenderedPageBreak/><w:t>Viết chương trình tính tổ</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>ng S(n) = x</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>^</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> + x^</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">4 + </w:t></w:r><w:r w:rsidRPr="001151B4"><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>⋯</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> + x^</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>2n</w:t></w:r><w:r w:rsidRPr="00C5583D"><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>, sử dụng PHP với n là biến cho trước</w:t></w:r><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>, ^ là ký tự số mũ</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tongParaRange.InsertXML($tongParaXml)

# --- Step 3: add vi-VN language tagging across the "tinh tong 1..500" block ---
$langBlockStart = Get-ParagraphRangeByText $d "Hãy thực hiện"
$langBlockEndPara = Get-ParagraphRangeByText $d "Hiển thị cấu trúc sau"
$langBlockRange = $d.Range($langBlockStart.Start, $langBlockEndPara.End)
$langBlockXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009D45CA" w:rsidRDefault="00D81E9D" w:rsidP="009D45CA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="47"/></w:numPr><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Hãy thực hiện tính tổng các số từ 1 tới 500</w:t></w:r><w:r w:rsidR="008553D2"><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> sử dụng PHP</w:t></w:r><w:r w:rsidR="00437282"><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p w:rsidR="00437282" w:rsidRDefault="00437282" w:rsidP="00437282"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Kết quả hiển thị sẽ có dạng:</w:t></w:r></w:p><w:p w:rsidR="00437282" w:rsidRDefault="00437282" w:rsidP="00437282"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/><w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r w:rsidRPr="00437282"><w:rPr><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Tổng các số từ 1 đến 500 = &lt;giá-trị&gt;</w:t></w:r></w:p><w:p w:rsidR="008553D2" w:rsidRDefault="008553D2" w:rsidP="008553D2"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="47"/></w:numPr><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Hiển thị cấu trúc sau sử dụng PHP</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$langBlockRange.InsertXML($langBlockXml)

Write-Output "edit complete"
